# Apply cryptos list update (2024-11-25 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to Text
# format first, otherwise Excel would silently coerce them to numbers
# (losing formatting like trailing zeros, e.g. "8.80" -> 8.8) instead of
# keeping them as the literal strings seen in the source data feed.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D13", "D16", "D18", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated coin rows (price / volume refreshed; a few rows also reordered
# with their neighbours as ranking shuffled)
$ws.Range("D2").Value = "97.798.52"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "3.375.39"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "252.31"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").Value = "658.19"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "1.45"
$ws.Range("E7").Value = "  -1.02%  "
$ws.Range("D8").Value = "0.422"
$ws.Range("E8").Value = "  -3.25%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "1.05"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("B10").Value = "USDC"
$ws.Range("C10").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "3.372.31"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("D13").Value = "41.54"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("D14").Value = "97.497.12"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("E15").Value = "  -5.66%  "
$ws.Range("D16").Value = "0.0000254"
$ws.Range("E16").Value = "  -4.85%  "
$ws.Range("D17").Value = "4.005.94"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "8.88"
$ws.Range("E18").Value = "  -2.88%  "
$ws.Range("D19").Value = "3.357.98"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").Value = "18.09"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").Value = "0.542"
$ws.Range("E21").Value = "  -9.31%  "
$ws.Range("D22").Value = "10.97"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "511.05"
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("E24").Value = "  -3.66%  "
$ws.Range("D25").Value = "7.09"
$ws.Range("E25").Value = "  +9.83%  "
$ws.Range("D26").Value = "0.0000199"
$ws.Range("E26").Value = "  -3.73%  "
$ws.Range("D27").Value = "96.58"
$ws.Range("E27").Value = "  -4.55%  "
$ws.Range("D28").Value = "12.33"
$ws.Range("E28").Value = "  -6.33%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.554.45"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "11.36"
$ws.Range("E30").Value = "  -5.80%  "
$ws.Range("E31").Value = "  -6.01%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").Value = "0.186"
$ws.Range("E33").Value = "  -6.35%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "2.57"
$ws.Range("E34").Value = "  +7.70%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "0.562"
$ws.Range("E36").Value = "  -3.87%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "28.70"
$ws.Range("E37").Value = "  -5.04%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "8.04"
$ws.Range("E38").Value = "  +1.28%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "1.50"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "529.24"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.152"
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "24.39"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").Value = "0.856"
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("D46").Value = "0.0426"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "2.30"
$ws.Range("E47").Value = "  +7.96%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "56.17"
$ws.Range("E48").Value = "  +9.02%  "
$ws.Range("B49").Value = "MantraDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D49").Value = "3.67"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "5.61"
$ws.Range("E50").Value = "  -6.72%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "8.58"
$ws.Range("E51").Value = "  -6.56%  "
